$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet: 展览 (sheet1)
$ws1.Range("F2").Value = 231
$ws1.Range("F3").Value = 1055
$ws1.Range("F4").Value = 526
$ws1.Range("F5").Value = 13857
$ws1.Range("F7").Value = 557
$ws1.Range("F8").Value = 4
$ws1.Range("F9").Value = 1780
$ws1.Range("F16").Value = 2
$ws1.Range("F17").Value = 13899
$ws1.Range("F19").Value = 621
$ws1.Range("F20").Value = 14959
$ws1.Range("F22").Value = 8243
$ws1.Range("F25").Value = 24
$ws1.Range("F28").Value = 163
$ws1.Range("F30").Value = 28
$ws1.Range("F31").Value = 1035
$ws1.Range("F32").Value = 7
$ws1.Range("F33").Value = 18
$ws1.Range("F36").Value = 6
$ws1.Range("F37").Value = 10
$ws1.Range("F38").Value = 218
$ws1.Range("F39").Value = 216
$ws1.Range("F40").Value = 390
$ws1.Range("F42").Value = 5078

# Sheet: 全部类型 (sheet4)
$ws4.Range("F2").Value = 231
$ws4.Range("F3").Value = 1055
$ws4.Range("F4").Value = 526
$ws4.Range("F5").Value = 13857
$ws4.Range("F7").Value = 557
$ws4.Range("F8").Value = 4
$ws4.Range("F9").Value = 1780
$ws4.Range("F16").Value = 2
$ws4.Range("F17").Value = 13899
$ws4.Range("F19").Value = 621
$ws4.Range("F20").Value = 14959
$ws4.Range("F22").Value = 8243
$ws4.Range("F25").Value = 24
$ws4.Range("F28").Value = 163
$ws4.Range("F30").Value = 28
$ws4.Range("F31").Value = 1035
$ws4.Range("F32").Value = 7
$ws4.Range("F33").Value = 18
$ws4.Range("F38").Value = 6
$ws4.Range("F39").Value = 10
$ws4.Range("F40").Value = 218
$ws4.Range("F41").Value = 216
$ws4.Range("F42").Value = 390
$ws4.Range("F44").Value = 5078
